# ---------------------------------------------------------------------------
# Applies the "Fixed update to excel issue" commit:
#   1. Rename "Weekly Quantity"!B1 header  "Requested quantity" -> "Weekly_PO_Qty"
#   2. Rename "Monthly Trend"!B1 header    "Requested quantity" -> "Monthly_PO_Qty"
#   3. Add a new "PO Forecast" worksheet (3rd tab) with a ds / PO_Forecast /
#      yhat_lower / yhat_upper forecast table.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. Weekly Quantity sheet: rename B1 header ----------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

# --- 2. Monthly Trend sheet: rename B1 header -------------------------------
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 3. Add the "PO Forecast" sheet as the 3rd tab --------------------------
$wsForecast = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

$headerRange = $wsForecast.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Data rows: ds, PO_Forecast, yhat_lower, yhat_upper
$data = @(
    @(45200.99999999999, 85, -269.9121718468795, 450.7108968157847),
    @(45214.99999999999, 106, -253.0073152950982, 482.5497989747831),
    @(45221.99999999999, 117, -234.8178634005335, 481.837404987615),
    @(45242.99999999999, 148, -230.9312956558949, 481.2698246674684),
    @(45298.99999999999, 233, -125.0656076614511, 623.9619983065046),
    @(45305.99999999999, 244, -122.3951231354813, 608.8868884022197),
    @(45312.99999999999, 254, -108.1521709970343, 631.5093184156981),
    @(45319.99999999999, 265, -102.4570312947978, 636.5197498656981),
    @(45326.99999999999, 275, -56.08417083677161, 630.7405549572914),
    @(45333.99999999999, 286, -79.10122397855098, 627.086816485435),
    @(45361.99999999999, 328, -31.56316329364366, 682.492502647866),
    @(45368.99999999999, 339, -27.61909139294445, 709.6859991285685),
    @(45375.99999999999, 349, 0.8777824878078072, 743.4685363059863),
    @(45389.99999999999, 371, 29.49809877743871, 734.4673187249846),
    @(45445.99999999999, 455, 94.76202146061573, 821.3871766766219),
    @(45452.99999999999, 466, 115.1225989911782, 818.9609212197529),
    @(45480.99999999999, 508, 157.3434598926093, 885.2568304766758),
    @(45508.99999999999, 551, 193.8747047985334, 911.0079423146054),
    @(45515.99999999999, 561, 215.1902824614369, 927.7466040305758),
    @(45522.99999999999, 572, 220.317990754272, 943.0646026983825),
    @(45557.99999999999, 625, 259.4760838903588, 980.8774764163805),
    @(45564.99999999999, 635, 274.735116415907, 1007.044943942112),
    @(45571.99999999999, 646, 280.185304473681, 993.7932636154964),
    @(45578.99999999999, 656, 311.3513437227781, 1028.930576661628),
    @(45585.99999999999, 667, 324.4456211426951, 1034.998561968919),
    @(45592.99999999999, 678, 326.33982017944, 1031.818909447295),
    @(45599.99999999999, 688, 314.4115951749717, 1039.255679124455),
    @(45606.99999999999, 699, 369.7012751559806, 1068.206674666646),
    @(45613.99999999999, 709, 349.3407671356542, 1089.914721680137),
    @(45620.99999999999, 720, 357.8034280265288, 1090.33440384236),
    @(45627.99999999999, 731, 377.6032386948386, 1089.662797263761)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $record = $data[$i]

    $dsCell = $wsForecast.Cells.Item($row, 1)
    $dsCell.Value = $record[0]
    $dsCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $wsForecast.Cells.Item($row, 2).Value = $record[1]
    $wsForecast.Cells.Item($row, 3).Value = $record[2]
    $wsForecast.Cells.Item($row, 4).Value = $record[3]
}
